$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.465.49"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.349.97"
$ws.Range("E3").Value = "  +5.89%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.17%  "
$ws.Range("D16").Value = "2.704.38"
$ws.Range("E16").Value = "  +6.13%  "
$ws.Range("D17").Value = "2.424.22"
$ws.Range("E17").Value = "  +8.60%  "
$ws.Range("D18").Value = "43.428.70"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "255.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.93%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0930"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.79%  "
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0377"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +11.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  +14.54%  "
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("E47").Value = "  +10.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.62%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.461"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.15%  "
